$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.459.65"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.888.56"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'243.60"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "'0.2896"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "'0.06489"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'22.18"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").Value = "'0.07755"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.887.33"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'95.65"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'0.7253"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "'5.190"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "'281.97"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "30.450.48"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'0.000007471"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "2.136.61"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'5.274"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'6.265"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").Value = "'163.83"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'9.080"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'18.86"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'1.891"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.333"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.09700"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").Value = "'4.271"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'4.143"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "'0.04859"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.01886"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").Value = "'2.814"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").Value = "'75.40"
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("D41").Value = "'6.220"
$ws.Range("D42").Value = "'1.983"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'0.4262"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.8255"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").Value = "'101.31"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "'9.603"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").Value = "'6.954"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "'35.13"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'909.47"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "'0.05749"
$ws.Range("E51").Value = "  +1.67%  "
